$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (numeric-looking percentages/prices) to be treated as plain
# text so Excel does not auto-convert them to numbers/percentages on assignment,
# matching the original inlineStr text cells in the workbook.
$ws.Range("D2:E50").NumberFormat = "@"

$ws.Range("E2").Value = "1.19%"
$ws.Range("D3").Value = "30.09"
$ws.Range("E3").Value = "11.61%"
$ws.Range("D4").Value = "5.162"
$ws.Range("E4").Value = "0.06%"
$ws.Range("D5").Value = "0.05742"
$ws.Range("E5").Value = "2.26%"
$ws.Range("D6").Value = "6.590"
$ws.Range("E6").Value = "1.44%"
$ws.Range("D7").Value = "0.8560"
$ws.Range("E7").Value = "4.68%"
$ws.Range("D8").Value = "0.8826"
$ws.Range("E8").Value = "6.09%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "0.1367"
$ws.Range("E9").Value = "2.74%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.03311"
$ws.Range("E10").Value = "7.73%"
$ws.Range("D11").Value = "0.06990"
$ws.Range("E11").Value = "0.98%"
$ws.Range("D12").Value = "0.02925"
$ws.Range("E12").Value = "0.92%"
$ws.Range("D13").Value = "0.09387"
$ws.Range("E13").Value = "0.08%"
$ws.Range("D14").Value = "0.001524"
$ws.Range("E14").Value = "0.71%"
$ws.Range("D15").Value = "0.04164"
$ws.Range("E15").Value = "-9.21%"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "0.0006018"
$ws.Range("E16").Value = "0.12%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "0.006149"
$ws.Range("E17").Value = "-0.27%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "3.509"
$ws.Range("E18").Value = "-3.71%"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "3.094"
$ws.Range("E19").Value = "2.35%"
$ws.Range("B20").Value = "BTSEToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D20").Value = "2.184"
$ws.Range("E20").Value = "-5.10%"
$ws.Range("B21").Value = "BitpandaEcosystemToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D21").Value = "0.3147"
$ws.Range("E21").Value = "1.13%"
$ws.Range("D22").Value = "0.1305"
$ws.Range("E22").Value = "1.07%"
$ws.Range("D23").Value = "3.606"
$ws.Range("E23").Value = "-3.70%"
$ws.Range("E24").Value = "2.62%"
$ws.Range("D25").Value = "0.001214"
$ws.Range("E25").Value = "-0.97%"
$ws.Range("D26").Value = "0.004501"
$ws.Range("E26").Value = "0.27%"
$ws.Range("D27").Value = "0.0001178"
$ws.Range("E27").Value = "20.16%"
$ws.Range("E28").Value = "-1.57%"
$ws.Range("D40").Value = "0.03787"
$ws.Range("E40").Value = "4.09%"
$ws.Range("D41").Value = "0.005749"
$ws.Range("E41").Value = "-5.46%"
$ws.Range("D42").Value = "0.1071"
$ws.Range("E42").Value = "1.93%"
$ws.Range("D43").Value = "0.002555"
$ws.Range("E43").Value = "-2.05%"
$ws.Range("D44").Value = "0.009984"
$ws.Range("E44").Value = "22.82%"
$ws.Range("D45").Value = "0.00005090"
$ws.Range("E45").Value = "-3.87%"
$ws.Range("E46").Value = "-0.22%"
$ws.Range("D47").Value = "0.08882"
$ws.Range("E47").Value = "-18.52%"
$ws.Range("D48").Value = "0.002720"
$ws.Range("E48").Value = "4.70%"
$ws.Range("D49").Value = "0.00002096"
$ws.Range("E49").Value = "-0.22%"
$ws.Range("D50").Value = "0.0001996"

# Restore default (unstyled) cell style for the D:E range so the saved XML does not
# retain an explicit number format style on cells that originally had none.
$ws.Range("D2:E50").Style = "Normal"

